# PfItDetailAdjust.xlsx — the "DATE" column-type label used by the two
# audit-trail columns (CreateDate / LastUpdate rows) is renamed to
# "TIMESTAMP" on the DBD sheet, and the sheet's active-cell selection is
# moved to D22 (the last cell touched).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$ws.Activate()

$ws.Range("D20").Value = "TIMESTAMP"
$ws.Range("D22").Value = "TIMESTAMP"

# Update the active cell selection to match the saved state.
$ws.Range("D22").Select()
